# Add a new submission row (row 13) to the overview table, for the RF
# submission that averages over tiles before predicting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new row's cell values. The order in which new (not-yet-seen)
# strings are first written determines their position in the shared
# string table, so write them in the same order as the target workbook.
$ws.Range("F13").Value = "weakly supervision with cv centers & average before predictions"
$ws.Range("A13").Value = "2023-02-27-1906_RF_centers_agg_pred.csv"
$ws.Range("G13").Value = "0.609 (0.016)"
$ws.Range("H13").Value = "Feb. 27, 2023, 6:20 p.m."
$ws.Range("B13").Value = "RandomForest"
$ws.Range("C13").Value = "MoCov"
$ws.Range("D13").Value = "-"
$ws.Range("E13").Value = "1 x 3"
$ws.Range("I13").Value = 0.605

# Grow the table ("Tabelle1") so it covers the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:I13"))

# Update the window selection to match where the user ended up after
# entering the new row.
$ws.Range("I14").Select()
